# Hoang tich hop chuc nang QuanLyLoaiDiaDiem
# - Fill in the "% done" / result columns (D9, E9) for the
#   "3.7 Quan ly thong tin Loai Dia Diem" row on the PhanCongNganHan sheet.
# - D9 gets a percentage number format applied (new cell style).
# - Update the sheet's active selection to E10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PhanCongNganHan")

$ws.Range("D9").Value2 = "75% (Làm được 3/4 chức năng)"
$ws.Range("E9").Value2 = "100% (10/06/2010)"

$ws.Range("D9").NumberFormat = "0%"

$ws.Range("E10").Select() | Out-Null
